$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value2 = 36.3798942565918
$ws.Cells.Item(2, 2).Value2 = 0.2311378717422485
$ws.Cells.Item(2, 3).Value2 = 33.12029266357422
$ws.Cells.Item(2, 4).Value2 = 34.84400939941406
$ws.Cells.Item(2, 5).Value2 = 38.1926383972168
$ws.Cells.Item(2, 6).Value2 = 5.072345733642578
$ws.Cells.Item(2, 7).Value2 = 132.3363342285156
$ws.Cells.Item(2, 8).Value2 = 66.71249389648438
$ws.Cells.Item(2, 9).Value2 = 320.8468017578125
$ws.Cells.Item(2, 10).Value2 = 247.9119415283203
$ws.Cells.Item(2, 11).Value2 = 0.7433949708938599
$ws.Cells.Item(2, 12).Value2 = 1.172901153564453
$ws.Cells.Item(2, 13).Value2 = 0.06797083467245102
$ws.Cells.Item(2, 14).Value2 = 0.5204110741615295
$ws.Cells.Item(2, 15).Value2 = 1.092594742774963
$ws.Cells.Item(2, 16).Value2 = 1.024623870849609
$ws.Cells.Item(2, 17).Value2 = 11.99386215209961
$ws.Cells.Item(2, 18).Value2 = 7.117724418640137
$ws.Cells.Item(2, 19).Value2 = 12.80570602416992
$ws.Cells.Item(2, 20).Value2 = 10.97202491760254
$ws.Cells.Item(2, 21).Value2 = 0.5354536771774292
$ws.Cells.Item(2, 22).Value2 = 1.571723461151123
$ws.Cells.Item(2, 23).Value2 = 10.09767913818359
$ws.Cells.Item(2, 24).Value2 = 2.589601516723633
$ws.Cells.Item(2, 25).Value2 = 15.44327926635742
$ws.Cells.Item(2, 26).Value2 = 0.9973151683807373
$ws.Cells.Item(2, 27).Value2 = 13.28406715393066
$ws.Cells.Item(2, 28).Value2 = 1.637460231781006
$ws.Cells.Item(2, 29).Value2 = -2.899353742599487
$ws.Cells.Item(2, 30).Value2 = -7.929656505584717
$ws.Cells.Item(2, 31).Value2 = 0.04105398058891296
$ws.Cells.Item(2, 32).Value2 = 1.054503083229065
$ws.Cells.Item(2, 33).Value2 = 1.405201315879822
$ws.Cells.Item(2, 34).Value2 = 0.7153990864753723
$ws.Cells.Item(2, 35).Value2 = 4.853835105895996
$ws.Cells.Item(2, 36).Value2 = 1.296501636505127
$ws.Cells.Item(2, 37).Value2 = 5.699641227722168
$ws.Cells.Item(2, 38).Value2 = 2.666948080062866
$ws.Cells.Item(2, 39).Value2 = 22.63392448425293
$ws.Cells.Item(2, 40).Value2 = 0.9566462635993958
$ws.Cells.Item(2, 41).Value2 = 717.261474609375
$ws.Cells.Item(2, 42).Value2 = 0.5391760468482971
$ws.Cells.Item(2, 43).Value2 = 1058.71533203125
$ws.Cells.Item(2, 44).Value2 = 0.3168524503707886
$ws.Cells.Item(2, 45).Value2 = -102.5290069580078
$ws.Cells.Item(2, 46).Value2 = -0.9027989506721497
$ws.Cells.Item(2, 47).Value2 = 1884.626220703125
$ws.Cells.Item(2, 48).Value2 = 0.2096151262521744
$ws.Cells.Item(2, 49).Value2 = 809.751220703125
$ws.Cells.Item(2, 50).Value2 = 0.4383590519428253
$ws.Cells.Item(2, 51).Value2 = -102.8067932128906
$ws.Cells.Item(2, 52).Value2 = -0.8393402695655823
$ws.Cells.Item(2, 53).Value2 = 2997.9375
$ws.Cells.Item(2, 54).Value2 = 0.1051927357912064
$ws.Cells.Item(2, 55).Value2 = 844.0130004882812
$ws.Cells.Item(2, 56).Value2 = 0.5375494360923767
$ws.Cells.Item(2, 57).Value2 = -104.4080657958984
$ws.Cells.Item(2, 58).Value2 = -0.8150216341018677
$ws.Cells.Item(2, 59).Value2 = -12.93214416503906
$ws.Cells.Item(2, 60).Value2 = -1.410694599151611
$ws.Cells.Item(2, 61).Value2 = 21.38112831115723
$ws.Cells.Item(2, 62).Value2 = 0.7761311531066895
$ws.Cells.Item(2, 63).Value2 = 0.08713994175195694
$ws.Cells.Item(2, 64).Value2 = 0.3648333251476288
$ws.Cells.Item(2, 65).Value2 = -0.02833130210638046
$ws.Cells.Item(2, 66).Value2 = -0.6479204893112183
$ws.Cells.Item(2, 67).Value2 = 0.9177460670471191
$ws.Cells.Item(2, 68).Value2 = 1.076560854911804
$ws.Cells.Item(2, 69).Value2 = 20.10033226013184
$ws.Cells.Item(2, 70).Value2 = 1.327543377876282
$ws.Cells.Item(2, 71).Value2 = 13.85799407958984
$ws.Cells.Item(2, 72).Value2 = 1.125352144241333
$ws.Cells.Item(2, 73).Value2 = 20.47581100463867
$ws.Cells.Item(2, 74).Value2 = 1.108245730400085
$ws.Cells.Item(2, 75).Value2 = -15.63738536834717
$ws.Cells.Item(2, 76).Value2 = -1.402887344360352
$ws.Cells.Item(2, 77).Value2 = -2.609011173248291
$ws.Cells.Item(2, 78).Value2 = 9.712672233581543
$ws.Cells.Item(2, 79).Value2 = 0.06369610875844955
$ws.Cells.Item(2, 80).Value2 = -0.002465426689013839
$ws.Cells.Item(2, 81).Value2 = 0.09334204345941544
$ws.Cells.Item(2, 82).Value2 = 2.728731870651245
$ws.Cells.Item(2, 83).Value2 = 2.427184581756592
$ws.Cells.Item(2, 84).Value2 = 0.2139999866485596
$ws.Cells.Item(2, 85).Value2 = 0.2148425579071045
$ws.Cells.Item(2, 86).Value2 = 0.18857142329216
$ws.Cells.Item(2, 87).Value2 = 0.2943238317966461
$ws.Cells.Item(2, 88).Value2 = -21.38953590393066
